$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 10
$ws_ALC.Range("H10").Value = 4666.3335
$ws_ALC.Range("I10").Value = 500
$ws_ALC.Range("J10").Value = 6749.5
$ws_ALC.Range("K10").Value = 500
$ws_ALC.Range("L10").Value = 6749.5
$ws_ALC.Range("M10").Value = -207
$ws_ALC.Range("N10").Value = -7335.5

# ALC row 18
$ws_ALC.Range("H18").Value = 900
$ws_ALC.Range("I18").Value = 900
$ws_ALC.Range("K18").Value = 900
$ws_ALC.Range("M18").Value = -616

# ALC row 70
$ws_ALC.Range("H70").Value = 2000
$ws_ALC.Range("I70").Value = 1400
$ws_ALC.Range("J70").Value = 3200
$ws_ALC.Range("K70").Value = 4200
$ws_ALC.Range("L70").Value = 9600
$ws_ALC.Range("M70").Value = -3930
$ws_ALC.Range("N70").Value = -10140

# ALC row 73
$ws_ALC.Range("H73").Value = 2000
$ws_ALC.Range("I73").Value = 1400
$ws_ALC.Range("J73").Value = 3200
$ws_ALC.Range("K73").Value = 4200
$ws_ALC.Range("L73").Value = 9600
$ws_ALC.Range("M73").Value = -3264
$ws_ALC.Range("N73").Value = -11472

# ALC row 141
$ws_ALC.Range("H141").Value = 7321.2856
$ws_ALC.Range("I141").Value = 7321.2856
$ws_ALC.Range("K141").Value = 21963.8568
$ws_ALC.Range("M141").Value = -16783.8568

# ARM row 19
$ws_ARM.Range("H19").Value = 5000
$ws_ARM.Range("I19").Value = 0
$ws_ARM.Range("J19").Value = 5000
$ws_ARM.Range("K19").Value = 0
$ws_ARM.Range("L19").Value = 5000
$ws_ARM.Range("M19").ClearContents()
$ws_ARM.Range("N19").Value = -5458

# ARM row 45
$ws_ARM.Range("H45").Value = 4315.8
$ws_ARM.Range("I45").Value = 1810
$ws_ARM.Range("K45").Value = 1810
$ws_ARM.Range("M45").Value = -1433

# ARM row 119
$ws_ARM.Range("H119").Value = 73699.5
$ws_ARM.Range("J119").Value = 73699.5
$ws_ARM.Range("L119").Value = 73699.5
$ws_ARM.Range("N119").Value = -83375.5

# ARM row 135
$ws_ARM.Range("H135").Value = 0
$ws_ARM.Range("J135").Value = 0
$ws_ARM.Range("L135").Value = 0
$ws_ARM.Range("N135").ClearContents()

# BSM row 40
$ws_BSM.Range("H40").Value = 19000
$ws_BSM.Range("I40").Value = 19000
$ws_BSM.Range("K40").Value = 19000
$ws_BSM.Range("M40").Value = -18735

# BSM row 80
$ws_BSM.Range("H80").Value = 530.5
$ws_BSM.Range("I80").Value = 682.6667
$ws_BSM.Range("J80").Value = 439.2
$ws_BSM.Range("K80").Value = 682.6667
$ws_BSM.Range("L80").Value = 439.2
$ws_BSM.Range("M80").Value = 315.3333
$ws_BSM.Range("N80").Value = -2435.2

# BSM row 83
$ws_BSM.Range("H83").Value = 530.5
$ws_BSM.Range("I83").Value = 682.6667
$ws_BSM.Range("J83").Value = 439.2
$ws_BSM.Range("K83").Value = 3413.3335
$ws_BSM.Range("L83").Value = 2196
$ws_BSM.Range("M83").Value = 1578.6665
$ws_BSM.Range("N83").Value = -12180

# BSM row 105
$ws_BSM.Range("H105").Value = 3655
$ws_BSM.Range("I105").Value = 3255.7693
$ws_BSM.Range("K105").Value = 3255.7693
$ws_BSM.Range("M105").Value = -1508.7693

# CUL row 16
$ws_CUL.Range("H16").Value = 2
$ws_CUL.Range("I16").Value = 2
$ws_CUL.Range("K16").Value = 6
$ws_CUL.Range("M16").Value = 167

# CUL row 19
$ws_CUL.Range("H19").Value = 9997
$ws_CUL.Range("J19").Value = 9997
$ws_CUL.Range("L19").Value = 29991
$ws_CUL.Range("N19").Value = -30339

# GSM row 20
$ws_GSM.Range("H20").Value = 51388.375
$ws_GSM.Range("J20").Value = 51388.375
$ws_GSM.Range("L20").Value = 51388.375
$ws_GSM.Range("N20").Value = -51878.375

# GSM row 80
$ws_GSM.Range("H80").Value = 950
$ws_GSM.Range("I80").Value = 900
$ws_GSM.Range("J80").Value = 1000
$ws_GSM.Range("K80").Value = 900
$ws_GSM.Range("L80").Value = 1000
$ws_GSM.Range("M80").Value = 98
$ws_GSM.Range("N80").Value = -2996

# GSM row 83
$ws_GSM.Range("H83").Value = 950
$ws_GSM.Range("I83").Value = 900
$ws_GSM.Range("J83").Value = 1000
$ws_GSM.Range("K83").Value = 4500
$ws_GSM.Range("L83").Value = 5000
$ws_GSM.Range("M83").Value = 492
$ws_GSM.Range("N83").Value = -14984

# GSM row 103
$ws_GSM.Range("H103").Value = 0
$ws_GSM.Range("J103").Value = 0
$ws_GSM.Range("L103").Value = 0
$ws_GSM.Range("N103").ClearContents()

# GSM row 107
$ws_GSM.Range("H107").Value = 651.625
$ws_GSM.Range("I107").Value = 474.66666
$ws_GSM.Range("K107").Value = 474.66666
$ws_GSM.Range("M107").Value = 1445.33334

# GSM row 122
$ws_GSM.Range("H122").Value = 3000
$ws_GSM.Range("I122").Value = 1666.6666
$ws_GSM.Range("K122").Value = 4999.9998
$ws_GSM.Range("M122").Value = -2549.9998

# GSM row 136
$ws_GSM.Range("H136").Value = 22730.4
$ws_GSM.Range("J136").Value = 27884
$ws_GSM.Range("L136").Value = 83652
$ws_GSM.Range("N136").Value = -88752

# LTW row 18
$ws_LTW.Range("H18").Value = 15000
$ws_LTW.Range("J18").Value = 15000
$ws_LTW.Range("L18").Value = 15000
$ws_LTW.Range("N18").Value = -15344

# LTW row 46
$ws_LTW.Range("H46").Value = 3242.6
$ws_LTW.Range("I46").Value = 2500
$ws_LTW.Range("J46").Value = 3985.2
$ws_LTW.Range("K46").Value = 2500
$ws_LTW.Range("L46").Value = 3985.2
$ws_LTW.Range("M46").Value = -2312
$ws_LTW.Range("N46").Value = -4361.2

# LTW row 68
$ws_LTW.Range("H68").Value = 4000.6667
$ws_LTW.Range("I68").Value = 0
$ws_LTW.Range("J68").Value = 4000.6667
$ws_LTW.Range("K68").Value = 0
$ws_LTW.Range("L68").Value = 4000.6667
$ws_LTW.Range("M68").ClearContents()
$ws_LTW.Range("N68").Value = -5498.6667

# LTW row 71
$ws_LTW.Range("H71").Value = 4000.6667
$ws_LTW.Range("I71").Value = 0
$ws_LTW.Range("J71").Value = 4000.6667
$ws_LTW.Range("K71").Value = 0
$ws_LTW.Range("L71").Value = 20003.3335
$ws_LTW.Range("M71").ClearContents()
$ws_LTW.Range("N71").Value = -27491.3335

# LTW row 82
$ws_LTW.Range("H82").Value = 2842.7856
$ws_LTW.Range("I82").Value = 968
$ws_LTW.Range("J82").Value = 3354.0908
$ws_LTW.Range("K82").Value = 968
$ws_LTW.Range("L82").Value = 3354.0908
$ws_LTW.Range("M82").Value = -607
$ws_LTW.Range("N82").Value = -4076.0908

# LTW row 85
$ws_LTW.Range("H85").Value = 2842.7856
$ws_LTW.Range("I85").Value = 968
$ws_LTW.Range("J85").Value = 3354.0908
$ws_LTW.Range("K85").Value = 968
$ws_LTW.Range("L85").Value = 3354.0908
$ws_LTW.Range("M85").Value = 280
$ws_LTW.Range("N85").Value = -5850.0908

# LTW row 99
$ws_LTW.Range("H99").Value = 90258.5
$ws_LTW.Range("I99").Value = 90258.5
$ws_LTW.Range("K99").Value = 90258.5
$ws_LTW.Range("M99").Value = -87263.5

# LTW row 119
$ws_LTW.Range("H119").Value = 20000
$ws_LTW.Range("J119").Value = 20000
$ws_LTW.Range("L119").Value = 20000
$ws_LTW.Range("N119").Value = -29676

# LTW row 132
$ws_LTW.Range("H132").Value = 7594.095
$ws_LTW.Range("I132").Value = 6375.4443
$ws_LTW.Range("J132").Value = 8508.083000000001
$ws_LTW.Range("K132").Value = 19126.3329
$ws_LTW.Range("L132").Value = 25524.249
$ws_LTW.Range("M132").Value = -16596.3329
$ws_LTW.Range("N132").Value = -30584.249

# LTW row 133
$ws_LTW.Range("H133").Value = 21000
$ws_LTW.Range("J133").Value = 21000
$ws_LTW.Range("L133").Value = 21000
$ws_LTW.Range("N133").Value = -26060

# WVR row 62
$ws_WVR.Range("H62").Value = 4822.727
$ws_WVR.Range("I62").Value = 4375
$ws_WVR.Range("J62").Value = 5360
$ws_WVR.Range("K62").Value = 4375
$ws_WVR.Range("L62").Value = 5360
$ws_WVR.Range("M62").Value = -3751
$ws_WVR.Range("N62").Value = -6608

# WVR row 65
$ws_WVR.Range("H65").Value = 4822.727
$ws_WVR.Range("I65").Value = 4375
$ws_WVR.Range("J65").Value = 5360
$ws_WVR.Range("K65").Value = 21875
$ws_WVR.Range("L65").Value = 26800
$ws_WVR.Range("M65").Value = -18755
$ws_WVR.Range("N65").Value = -33040

# WVR row 68
$ws_WVR.Range("H68").Value = 35000
$ws_WVR.Range("J68").Value = 35000
$ws_WVR.Range("L68").Value = 35000
$ws_WVR.Range("N68").Value = -36622

# WVR row 71
$ws_WVR.Range("H71").Value = 35000
$ws_WVR.Range("J71").Value = 35000
$ws_WVR.Range("L71").Value = 105000
$ws_WVR.Range("N71").Value = -113112

# WVR row 75
$ws_WVR.Range("H75").Value = 25000
$ws_WVR.Range("J75").Value = 25000
$ws_WVR.Range("L75").Value = 25000
$ws_WVR.Range("N75").Value = -26872

# WVR row 78
$ws_WVR.Range("H78").Value = 25000
$ws_WVR.Range("J78").Value = 25000
$ws_WVR.Range("L78").Value = 75000
$ws_WVR.Range("N78").Value = -84360

# WVR row 86
$ws_WVR.Range("H86").Value = 47555
$ws_WVR.Range("J86").Value = 47555
$ws_WVR.Range("L86").Value = 47555
$ws_WVR.Range("N86").Value = -49801

# WVR row 89
$ws_WVR.Range("H89").Value = 47555
$ws_WVR.Range("J89").Value = 47555
$ws_WVR.Range("L89").Value = 237775
$ws_WVR.Range("N89").Value = -249007

# WVR row 119
$ws_WVR.Range("H119").Value = 72473.25
$ws_WVR.Range("J119").Value = 72473.25
$ws_WVR.Range("L119").Value = 72473.25
$ws_WVR.Range("N119").Value = -82149.25

# WVR row 132
$ws_WVR.Range("H132").Value = 1533.3
$ws_WVR.Range("I132").Value = 1323.5
$ws_WVR.Range("K132").Value = 3970.5
$ws_WVR.Range("M132").Value = -1440.5

# WVR row 137
$ws_WVR.Range("H137").Value = 89460.5
$ws_WVR.Range("J137").Value = 89460.5
$ws_WVR.Range("L137").Value = 89460.5
$ws_WVR.Range("N137").Value = -99660.5
